# Update North Dakota_A transition-matrix probabilities.
# The simulation engine now accounts for more games and an updated
# (faster) simulate-game routine, which shifts several transition
# probabilities in the state matrix on Sheet1. Apply the refreshed
# values cell-by-cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1891891891891892
$ws.Range("C2").Value = 0.543918918918919
$ws.Range("J2").Value = 0.05912162162162162
$ws.Range("P2").Value = 0.1334459459459459
$ws.Range("S2").Value = 0.07432432432432433
$ws.Range("B3").Value = 0.005865102639296188
$ws.Range("C3").Value = 0.02932551319648094
$ws.Range("J3").Value = 0.08797653958944282
$ws.Range("P3").Value = 0.750733137829912
$ws.Range("S3").Value = 0.126099706744868
$ws.Range("J4").Value = 0.1265822784810127
$ws.Range("P4").Value = 0.5949367088607594
$ws.Range("S4").Value = 0.2784810126582278
$ws.Range("B6").Value = 0.07657657657657657
$ws.Range("D6").Value = 0.006756756756756757
$ws.Range("E6").Value = 0.002252252252252252
$ws.Range("F6").Value = 0.06756756756756757
$ws.Range("J6").Value = 0.3018018018018018
$ws.Range("O6").Value = 0.01576576576576576
$ws.Range("Q6").Value = 0.1644144144144144
$ws.Range("R6").Value = 0.07657657657657657
$ws.Range("S6").Value = 0.2882882882882883
$ws.Range("B7").Value = 0.1198910081743869
$ws.Range("D7").Value = 0.0108991825613079
$ws.Range("E7").Value = 0.002724795640326975
$ws.Range("F7").Value = 0.04632152588555858
$ws.Range("J7").Value = 0.1689373297002725
$ws.Range("O7").Value = 0.008174386920980926
$ws.Range("Q7").Value = 0.1907356948228883
$ws.Range("R7").Value = 0.1008174386920981
$ws.Range("S7").Value = 0.3514986376021799
$ws.Range("B8").Value = 0.09619450317124736
$ws.Range("D8").Value = 0.02219873150105708
$ws.Range("F8").Value = 0.08033826638477801
$ws.Range("J8").Value = 0.1691331923890063
$ws.Range("O8").Value = 0.01162790697674419
$ws.Range("Q8").Value = 0.1913319238900634
$ws.Range("R8").Value = 0.1088794926004228
$ws.Range("S8").Value = 0.3202959830866808
$ws.Range("B9").Value = 0.09302325581395349
$ws.Range("D9").Value = 0.02616279069767442
$ws.Range("F9").Value = 0.07848837209302326
$ws.Range("J9").Value = 0.1598837209302326
$ws.Range("O9").Value = 0.02034883720930233
$ws.Range("Q9").Value = 0.1540697674418605
$ws.Range("R9").Value = 0.1220930232558139
$ws.Range("S9").Value = 0.3459302325581395
$ws.Range("B10").Value = 0.09592592592592593
$ws.Range("D10").Value = 0.01629629629629629
$ws.Range("E10").Value = 0.0007407407407407407
$ws.Range("F10").Value = 0.06148148148148148
$ws.Range("J10").Value = 0.2455555555555556
$ws.Range("O10").Value = 0.02592592592592593
$ws.Range("Q10").Value = 0.2055555555555555
$ws.Range("R10").Value = 0.07777777777777778
$ws.Range("S10").Value = 0.2707407407407407
$ws.Range("G11").Value = 0.1382978723404255
$ws.Range("J11").Value = 0.07872340425531915
$ws.Range("K11").Value = 0.1787234042553192
$ws.Range("L11").Value = 0.5957446808510638
$ws.Range("S11").Value = 0.008510638297872341
$ws.Range("G12").Value = 0.7929824561403509
$ws.Range("J12").Value = 0.1508771929824561
$ws.Range("K12").Value = 0.01052631578947368
$ws.Range("L12").Value = 0.02456140350877193
$ws.Range("S12").Value = 0.02105263157894737
$ws.Range("F13").Value = 0.009523809523809525
$ws.Range("G13").Value = 0.7238095238095238
$ws.Range("J13").Value = 0.2476190476190476
$ws.Range("S13").Value = 0.01904761904761905
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.5
$ws.Range("S14").Value = 0.1666666666666667
$ws.Range("F15").Value = 0.02638522427440633
$ws.Range("H15").Value = 0.2005277044854881
$ws.Range("I15").Value = 0.0554089709762533
$ws.Range("J15").Value = 0.3693931398416886
$ws.Range("K15").Value = 0.0633245382585752
$ws.Range("M15").Value = 0.01846965699208443
$ws.Range("N15").Value = 0.002638522427440633
$ws.Range("O15").Value = 0.05277044854881267
$ws.Range("S15").Value = 0.2110817941952507
$ws.Range("F16").Value = 0.01612903225806452
$ws.Range("H16").Value = 0.1935483870967742
$ws.Range("I16").Value = 0.09677419354838709
$ws.Range("J16").Value = 0.4327956989247312
$ws.Range("K16").Value = 0.08870967741935484
$ws.Range("M16").Value = 0.03225806451612903
$ws.Range("N16").Value = 0.002688172043010753
$ws.Range("O16").Value = 0.05913978494623656
$ws.Range("S16").Value = 0.07795698924731183
$ws.Range("F17").Value = 0.01841820151679307
$ws.Range("H17").Value = 0.2275189599133261
$ws.Range("I17").Value = 0.08450704225352113
$ws.Range("J17").Value = 0.4366197183098591
$ws.Range("K17").Value = 0.0866738894907909
$ws.Range("M17").Value = 0.02058504875406284
$ws.Range("N17").Value = 0.001083423618634886
$ws.Range("O17").Value = 0.05850487540628386
$ws.Range("S17").Value = 0.06608884073672806
$ws.Range("F18").Value = 0.02830188679245283
$ws.Range("H18").Value = 0.1839622641509434
$ws.Range("I18").Value = 0.06839622641509434
$ws.Range("J18").Value = 0.4834905660377358
$ws.Range("K18").Value = 0.09905660377358491
$ws.Range("M18").Value = 0.01886792452830189
$ws.Range("N18").Value = 0.002358490566037736
$ws.Range("O18").Value = 0.0589622641509434
$ws.Range("S18").Value = 0.05660377358490566
$ws.Range("F19").Value = 0.01863684771033014
$ws.Range("H19").Value = 0.2300319488817891
$ws.Range("I19").Value = 0.08093716719914804
$ws.Range("J19").Value = 0.3940362087326943
$ws.Range("K19").Value = 0.09531416400425985
$ws.Range("M19").Value = 0.02928647497337593
$ws.Range("N19").Value = 0.001064962726304579
$ws.Range("O19").Value = 0.05644302449414271
$ws.Range("S19").Value = 0.09424920127795527
